$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Sala de Aula"

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Nota dos alunos"
